# Update countries & provincias Spain
# Refreshes the COVID country/province stats table with newer data,
# which also changes the ranking order for a couple of tied/near-tied
# countries (Argentina/Banglades and Montserrat/Islas Malvinas).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last refreshed" timestamp shown at the top of the sheet
$ws.Range("A1").Value = "Datos actualizados a 13 de Agosto de 2020 a las 12:26"

# Estados Unidos (row 4) - new case counts
$ws.Range("B4").Value = 5360488
$ws.Range("C4").Value = 186
$ws.Range("D4").Value = 2813125
$ws.Range("E4").Value = 2378231
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 169132

# Row 14 - new case counts
$ws.Range("B14").Value = 336324
$ws.Range("C14").Value = 2625
$ws.Range("D14").Value = 292058
$ws.Range("E14").Value = 25104
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 174
$ws.Range("H14").Value = 19162

# Row 18 now becomes Banglades (overtook Argentina in ranking)
$ws.Range("A18").Value = "Banglades"
$ws.Range("B18").Value = 269115
$ws.Range("C18").Value = 2617
$ws.Range("D18").Value = 154871
$ws.Range("E18").Value = 110687
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 44
$ws.Range("H18").Value = 3557

# Row 19 now becomes Argentina (data unchanged from previous refresh)
$ws.Range("A19").Value = "Argentina"
$ws.Range("B19").Value = 268574
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 187283
$ws.Range("E19").Value = 76078
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 5213

# Row 43 - new case counts
$ws.Range("B43").Value = 66631
$ws.Range("C43").Value = 1454
$ws.Range("D43").Value = 31547
$ws.Range("E43").Value = 32224
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 53
$ws.Range("H43").Value = 2860

# Row 73 - new case counts
$ws.Range("B73").Value = 21993
$ws.Range("C73").Value = 349
$ws.Range("D73").Value = 10254
$ws.Range("E73").Value = 11155
$ws.Range("F73").Value = 0

# Row 88 - new case counts
$ws.Range("B88").Value = 9129
$ws.Range("C88").Value = 15
$ws.Range("D88").Value = 8821
$ws.Range("E88").Value = 183
$ws.Range("F88").Value = 0

# Row 213 now becomes Islas Malvinas
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214 now becomes Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
